$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# Filter strings changes, delete line if True for end clearing method.
$ws.Range("B16").Value = $true
$ws.Range("B24").Value = $false
